# Regenerate merged AHB files
# - Rename the "_old" / "_new" header suffixes to "_FV2310" / "_FV2404"
# - Turn the data range into a proper Excel Table (ListObject)
# - Freeze the header row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$baseNames = @("Segmentname", "Segmentgruppe", "Segment", "Datenelement", "Segment ID", "Code", "Qualifier", "Beschreibung", "Bedingungsausdruck", "Bedingung")
$oldCols   = @("A", "B", "C", "D", "E", "F", "G", "H", "I", "J")
$newCols   = @("L", "M", "N", "O", "P", "Q", "R", "S", "T", "U")

for ($i = 0; $i -lt $baseNames.Length; $i++) {
    $ws.Range($oldCols[$i] + "1").Value = $baseNames[$i] + "_FV2310"
    $ws.Range($newCols[$i] + "1").Value = $baseNames[$i] + "_FV2404"
}

# Convert the used range A1:U69 into an Excel Table with a header row
$tableRange = $ws.Range("A1:U69")
$lo = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $tableRange, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$lo.Name = "Table1"

# Freeze the header row (row 1)
$ws.Range("A2").Select() | Out-Null
($excel.ActiveWindow.FreezePanes = $true) | Out-Null

Write-Host "done"
